$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format first so numeric-looking strings
# (e.g. "1.007", "5.380") are stored as text, matching the source data,
# instead of being auto-parsed into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.552.26'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '1.844.26'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -1.24%  '
$ws.Range("D5").Value = '334.27'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").Value = '0.4634'
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("D8").Value = '0.3860'
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("D9").Value = '45.96'
$ws.Range("E9").Value = '  -3.15%  '
$ws.Range("D10").Value = '0.07897'
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("D11").Value = '0.9983'
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").Value = '21.47'
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.848.65'
$ws.Range("E13").Value = '  -1.95%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.965'
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").Value = '7.131'
$ws.Range("D16").Value = '1.007'
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("D17").Value = '88.41'
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").Value = '0.06685'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").Value = '0.00001035'
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("D20").Value = '17.13'
$ws.Range("E21").Value = '  -1.15%  '
$ws.Range("D22").Value = '27.537.28'
$ws.Range("E22").Value = '  -1.44%  '
$ws.Range("D23").Value = '5.380'
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").Value = '2.308'
$ws.Range("E25").Value = '  -1.91%  '
$ws.Range("D26").Value = '158.51'
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").Value = '2.059.07'
$ws.Range("E27").Value = '  -2.57%  '
$ws.Range("D28").Value = '19.50'
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("D29").Value = '2.112'
$ws.Range("E29").Value = '  +2.08%  '
$ws.Range("D30").Value = '5.411'
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").Value = '119.79'
$ws.Range("E31").Value = '  -1.04%  '
$ws.Range("D32").Value = '0.9765'
$ws.Range("E32").Value = '  +1.80%  '
$ws.Range("D33").Value = '0.09407'
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '3.601'
$ws.Range("E34").Value = '  -1.62%  '
$ws.Range("D35").Value = '5.300'
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = '1.330'
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("D37").Value = '0.06023'
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("D38").Value = '0.02224'
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("D39").Value = '8.303'
$ws.Range("E39").Value = '  +2.20%  '
$ws.Range("D40").Value = '1.179'
$ws.Range("E40").Value = '  -3.34%  '
$ws.Range("D41").Value = '0.5891'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").Value = '0.1861'
$ws.Range("E42").Value = '  -1.23%  '
$ws.Range("D43").Value = '10.33'
$ws.Range("E43").Value = '  +1.81%  '
$ws.Range("D44").Value = '1.238'
$ws.Range("E44").Value = '  -2.30%  '
$ws.Range("D45").Value = '0.5586'
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("D46").Value = '12.18'
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").Value = '1.909'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").Value = '0.06696'
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").Value = '110.19'
$ws.Range("E49").Value = '  -2.95%  '
$ws.Range("D50").Value = '1.049'
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").Value = '1.006'
$ws.Range("E51").Value = '  -1.16%  '

# Restore the default cell style on the Price column so no stray
# number-format style lingers on cells that did not need it.
$ws.Range("D2:D51").Style = "Normal"
